# Updates the weekly Hortaliza / Achicoria data:
#  - Adds a new data row (row 27) at the bottom of the table, duplicating the
#    row layout (columns A:R) used by the existing data rows.
#  - Re-populates the Fecha (D), Volumen (J), Precio minimo (K),
#    Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P)
#    columns for rows 2-27 with the refreshed weekly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by one row: copy the existing template row (row 2, which
# contains all the constant columns B/C/E/F/G/H/I/N/O/Q/R already shared by
# every data row) down into the new row 27.
$ws.Range("A2:R2").Copy($ws.Range("A27:R27"))
$excel.CutCopyMode = 0

# Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) for every data row.
$data = @(
    @(2, 44371, 34, 5500, 6000, 5750, 359),
    @(3, 44355, 25, 6000, 6000, 6000, 375),
    @(4, 44313, 34, 6000, 6000, 6000, 375),
    @(5, 44407, 45, 5500, 6000, 5744, 359),
    @(6, 44438, 34, 5000, 6000, 5500, 344),
    @(7, 44455, 52, 5000, 6000, 5500, 344),
    @(8, 44573, 34, 8000, 8000, 8000, 500),
    @(9, 44467, 52, 5000, 6000, 5500, 344),
    @(10, 44306, 50, 6000, 6000, 6000, 375),
    @(11, 44308, 70, 6000, 6000, 6000, 375),
    @(12, 44403, 43, 6000, 6000, 6000, 375),
    @(13, 44350, 25, 6000, 6000, 6000, 375),
    @(14, 44589, 52, 8000, 8000, 8000, 500),
    @(15, 44341, 51, 5500, 6000, 5755, 360),
    @(16, 44328, 160, 6000, 6000, 6000, 375),
    @(17, 44330, 120, 6000, 6000, 6000, 375),
    @(18, 44582, 52, 7000, 7000, 7000, 438),
    @(19, 44698, 34, 6000, 7000, 6500, 406),
    @(20, 44358, 52, 6000, 6000, 6000, 375),
    @(21, 44691, 61, 6000, 7000, 6508, 407),
    @(22, 44363, 160, 5500, 6000, 5750, 359),
    @(23, 44477, 25, 6000, 6000, 6000, 375),
    @(24, 44474, 52, 5000, 6000, 5500, 344),
    @(25, 44575, 61, 8000, 8000, 8000, 500),
    @(26, 44376, 43, 4500, 5000, 4756, 297),
    @(27, 44442, 25, 6000, 7000, 6480, 405)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]  # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]  # P: Precio $/Kg
}
